# "added cities and mission boards"
#
# 1. Bump the cached "today" date field (datetimeFigureOut) on the slide
#    master and every slide layout from 8/4/2020 -> 8/5/2020.
# 2. Slide 3 ("Planets"): append new city names to three existing lists
#    by adding runs right after the last run of the relevant paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders (slide master + all custom layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "8/5/2020"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# Helper: append `newText` right after an existing paragraph's text,
# then re-split the newly appended span into its own run(s) so the
# appended text doesn't just get absorbed into the previous run.
# `topRange` is the shape's full TextFrame.TextRange (used for absolute
# Characters() addressing, since Paragraph.Parent/.Index aren't usable
# here); `splits` is an array of substring lengths (e.g. @(2) to carve
# the first 2 new characters into their own run, leaving the remainder
# in a subsequent run).
# ---------------------------------------------------------------------
function Add-RunsAfterParagraph($topRange, $para, [string]$newText, [int[]]$splits) {
    $start = $para.Start
    $origLen = $para.Text.TrimEnd("`r").Length
    [void]$para.InsertAfter($newText)
    $pos = $start + $origLen
    foreach ($len in $splits) {
        $chunk = $topRange.Characters($pos, $len)
        $chunk.Text = $chunk.Text
        $pos += $len
    }
}

# ---------------------------------------------------------------------
# 2) Slide 3 ("Planets") city / mission-board additions
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

# "FederationPlanets" content placeholder -> Natotis paragraph gets " -Zhul"
$contentPH = $slide3.Shapes.Item(2)
$contentTR = $contentPH.TextFrame.TextRange
$natotisPara = $contentTR.Paragraphs(2)
Add-RunsAfterParagraph $contentTR $natotisPara " -Zhul" @(2)

# "TextBox 3" -> Otov paragraph gets "- Xigow"
$textBox3 = $slide3.Shapes.Item(3)
$tb3TR = $textBox3.TextFrame.TextRange
$otovPara = $tb3TR.Paragraphs(3)
Add-RunsAfterParagraph $tb3TR $otovPara "- Xigow" @(2)

# "TextBox 3" -> Larvis paragraph gets a trailing "-"
$larvisPara = $tb3TR.Paragraphs(4)
[void]$larvisPara.InsertAfter("-")
